$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title --------------------------------------------------------------
$ws.Range("B2").Value = 'Densidades de cada elemento'

# --- Header row -----------------------------------------------------------
$ws.Range("B3").Value = 'Elemento'
$ws.Range("C3").Value = 'Densidad [$T/m^3$]'

# --- Data rows --------------------------------------------------------
$ws.Range("B4").Value = '$\gamma_{Hormigon}$'
$ws.Range("C4").Value = 2.5

$ws.Range("B5").Value = '$\gamma_{Estuco}$'
$ws.Range("C5").Value = 2

$ws.Range("B6").Value = '$\gamma_{Sobrelosa}$'
$ws.Range("C6").Value = 1.5

# --- Formatting -----------------------------------------------------------
# Apply the thin box border to the whole table first so every cell shares
# the same border definition.
$table = $ws.Range("B3:C6")
$table.Borders.LineStyle = 1

# Header: bold + centered
$header = $ws.Range("B3:C3")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

# Row labels (left column of data rows): right aligned
$labels = $ws.Range("B4:B6")
$labels.HorizontalAlignment = -4152

# Values (right column of data rows): centered
$values = $ws.Range("C4:C6")
$values.HorizontalAlignment = -4108

# --- Column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.26
$ws.Columns.Item(2).ColumnWidth = 19.26
$ws.Columns.Item(3).ColumnWidth = 16.76
$ws.Columns.Item(4).ColumnWidth = 3.42

# --- Comment on C6 ------------------------------------------------------
$comment = $ws.Range("C6").AddComment("Autor:`nInvestigar")

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection (matches the cursor position the author left behind) ---
$ws.Range("E8").Select()
